# Update on 2018-04-30, 支出生活费400
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 37 (item #35) was blank; fill it in exactly like row 36 (item #34),
# which is the preceding "生活费" expense entry, then overwrite with the
# new period's values.
$ws.Range("C36:G36").Copy() | Out-Null
$ws.Range("C37:G37").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("C37").Value = "支出"
$ws.Range("D37").Value = 400
$ws.Range("E37").Value = "2018-04-30"
$ws.Range("F37").Value = "生活费"
$ws.Range("G37").Value = "生活费(5/1-5/10)"

# Match the scrolled/selected view state from the saved workbook.
$ws.Range("E38").Select() | Out-Null
